$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2368.3333
$ws.Range("I41").Value = 3282.1428
$ws.Range("J41").Value = 1568.75
$ws.Range("K41").Value = 3282.1428
$ws.Range("L41").Value = 1568.75
$ws.Range("M41").Value = -2842.1428
$ws.Range("N41").Value = -2448.75
$ws.Range("H76").Value = 3564.5833
$ws.Range("I76").Value = 3665.7896
$ws.Range("K76").Value = 3665.7896
$ws.Range("M76").Value = -3350.7896
$ws.Range("H79").Value = 3564.5833
$ws.Range("I79").Value = 3665.7896
$ws.Range("K79").Value = 3665.7896
$ws.Range("M79").Value = -2573.7896
$ws.Range("H86").Value = 7375
$ws.Range("I86").Value = 8666.666999999999
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 8666.666999999999
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -7543.666999999999
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 7375
$ws.Range("I89").Value = 8666.666999999999
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 43333.335
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -37717.335
$ws.Range("N89").Value = -28732
$ws.Range("H137").Value = 1141.5
$ws.Range("I137").Value = 811.1627999999999
$ws.Range("K137").Value = 2433.4884
$ws.Range("M137").Value = 116.5116000000003
$ws.Range("H138").Value = 1388.3738
$ws.Range("I138").Value = 808.871
$ws.Range("J138").Value = 1652.5588
$ws.Range("K138").Value = 2426.613
$ws.Range("L138").Value = 4957.6764
$ws.Range("M138").Value = 2713.387
$ws.Range("N138").Value = -15237.6764
$ws.Range("H141").Value = 904.44446
$ws.Range("I141").Value = 767.5
$ws.Range("K141").Value = 2302.5
$ws.Range("M141").Value = 2877.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4779.921
$ws.Range("I32").Value = 4318.6826
$ws.Range("J32").Value = 7015.154
$ws.Range("K32").Value = 4318.6826
$ws.Range("L32").Value = 7015.154
$ws.Range("M32").Value = -4031.6826
$ws.Range("N32").Value = -7589.154
$ws.Range("H61").Value = 76925040
$ws.Range("I61").Value = 125001310
$ws.Range("K61").Value = 125001310
$ws.Range("M61").Value = -125001098
$ws.Range("H74").Value = 1139.5588
$ws.Range("I74").Value = 926.80646
$ws.Range("J74").Value = 3338
$ws.Range("K74").Value = 926.80646
$ws.Range("L74").Value = 3338
$ws.Range("M74").Value = -52.80646000000002
$ws.Range("N74").Value = -5086
$ws.Range("H77").Value = 1139.5588
$ws.Range("I77").Value = 926.80646
$ws.Range("J77").Value = 3338
$ws.Range("K77").Value = 4634.0323
$ws.Range("L77").Value = 16690
$ws.Range("M77").Value = -266.0322999999999
$ws.Range("N77").Value = -25426
$ws.Range("H132").Value = 2216.4
$ws.Range("I132").Value = 2024.2941
$ws.Range("J132").Value = 2624.625
$ws.Range("K132").Value = 6072.8823
$ws.Range("L132").Value = 7873.875
$ws.Range("M132").Value = -3542.8823
$ws.Range("N132").Value = -12933.875
$ws.Range("H136").Value = 76925040
$ws.Range("I136").Value = 125001310
$ws.Range("K136").Value = 375003930
$ws.Range("M136").Value = -375001380

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3692.4634
$ws.Range("I134").Value = 1025.2258
$ws.Range("J134").Value = 11960.9
$ws.Range("K134").Value = 3075.6774
$ws.Range("L134").Value = 35882.7
$ws.Range("M134").Value = -540.6773999999996
$ws.Range("N134").Value = -40952.7

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 40001270
$ws.Range("I16").Value = 41667864
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 41667864
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -41667577
$ws.Range("N16").Value = -3574
$ws.Range("H31").Value = 2138.6428
$ws.Range("I31").Value = 2219.1
$ws.Range("J31").Value = 1937.5
$ws.Range("K31").Value = 2219.1
$ws.Range("L31").Value = 1937.5
$ws.Range("M31").Value = -1924.1
$ws.Range("N31").Value = -2527.5
$ws.Range("H34").Value = 2138.6428
$ws.Range("I34").Value = 2219.1
$ws.Range("J34").Value = 1937.5
$ws.Range("K34").Value = 2219.1
$ws.Range("L34").Value = 1937.5
$ws.Range("M34").Value = -2017.1
$ws.Range("N34").Value = -2341.5
$ws.Range("H58").Value = 1102.7179
$ws.Range("I58").Value = 1063.742
$ws.Range("J58").Value = 1253.75
$ws.Range("K58").Value = 1063.742
$ws.Range("L58").Value = 1253.75
$ws.Range("M58").Value = -860.742
$ws.Range("N58").Value = -1659.75
$ws.Range("H112").Value = 35100.285
$ws.Range("J112").Value = 38450.332
$ws.Range("L112").Value = 38450.332
$ws.Range("N112").Value = -41404.332
$ws.Range("H113").Value = 40001270
$ws.Range("I113").Value = 41667864
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 41667864
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -41665694
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 1971.8823
$ws.Range("I132").Value = 1623.826
$ws.Range("K132").Value = 4871.478
$ws.Range("M132").Value = -2341.478
$ws.Range("H134").Value = 18519714
$ws.Range("I134").Value = 1134.9412
$ws.Range("J134").Value = 50001300
$ws.Range("K134").Value = 3404.8236
$ws.Range("L134").Value = 150003900
$ws.Range("M134").Value = -869.8235999999997
$ws.Range("N134").Value = -150008970
$ws.Range("H136").Value = 1102.7179
$ws.Range("I136").Value = 1063.742
$ws.Range("J136").Value = 1253.75
$ws.Range("K136").Value = 3191.226
$ws.Range("L136").Value = 3761.25
$ws.Range("M136").Value = -641.2259999999997
$ws.Range("N136").Value = -8861.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 468407.12
$ws.Range("I4").Value = 50016.61
$ws.Range("J4").Value = 1544268.4
$ws.Range("K4").Value = 150049.83
$ws.Range("L4").Value = 4632805.199999999
$ws.Range("M4").Value = -149937.83
$ws.Range("N4").Value = -4633029.199999999
$ws.Range("H5").Value = 1815.1875
$ws.Range("I5").Value = 2074
$ws.Range("J5").Value = 1038.75
$ws.Range("K5").Value = 6222
$ws.Range("L5").Value = 3116.25
$ws.Range("M5").Value = -6110
$ws.Range("N5").Value = -3340.25
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 12
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 101
$ws.Range("H7").Value = 378.64285
$ws.Range("I7").Value = 455.44446
$ws.Range("J7").Value = 240.4
$ws.Range("K7").Value = 1366.33338
$ws.Range("L7").Value = 721.2
$ws.Range("M7").Value = -1254.33338
$ws.Range("N7").Value = -945.2
$ws.Range("H109").Value = 85800.664
$ws.Range("I109").Value = 167768
$ws.Range("J109").Value = 3833.3333
$ws.Range("K109").Value = 503304
$ws.Range("L109").Value = 11499.9999
$ws.Range("M109").Value = -502264
$ws.Range("N109").Value = -13579.9999
$ws.Range("H113").Value = 686.90625
$ws.Range("I113").Value = 498.25
$ws.Range("J113").Value = 713.8570999999999
$ws.Range("K113").Value = 1494.75
$ws.Range("L113").Value = 2141.5713
$ws.Range("M113").Value = 675.25
$ws.Range("N113").Value = -6481.5713
$ws.Range("H122").Value = 888.4483
$ws.Range("I122").Value = 693.0909
$ws.Range("J122").Value = 1007.8333
$ws.Range("K122").Value = 6237.8181
$ws.Range("L122").Value = 9070.4997
$ws.Range("M122").Value = -3787.8181
$ws.Range("N122").Value = -13970.4997
$ws.Range("H131").Value = 23259160
$ws.Range("J131").Value = 3723.9211
$ws.Range("L131").Value = 11171.7633
$ws.Range("N131").Value = -21251.7633
$ws.Range("H132").Value = 1190
$ws.Range("I132").Value = 1190
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10710
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8180
$ws.Range("H135").Value = 1815.1875
$ws.Range("I135").Value = 2074
$ws.Range("J135").Value = 1038.75
$ws.Range("K135").Value = 18666
$ws.Range("L135").Value = 9348.75
$ws.Range("M135").Value = -16131
$ws.Range("N135").Value = -14418.75
$ws.Range("N6").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50002956
$ws.Range("I70").Value = 41669684
$ws.Range("J70").Value = 66669500
$ws.Range("K70").Value = 41669684
$ws.Range("L70").Value = 66669500
$ws.Range("M70").Value = -41669414
$ws.Range("N70").Value = -66670040
$ws.Range("H73").Value = 50002956
$ws.Range("I73").Value = 41669684
$ws.Range("J73").Value = 66669500
$ws.Range("K73").Value = 41669684
$ws.Range("L73").Value = 66669500
$ws.Range("M73").Value = -41668748
$ws.Range("N73").Value = -66671372
$ws.Range("H132").Value = 3060.8
$ws.Range("I132").Value = 2694.3125
$ws.Range("J132").Value = 3712.3333
$ws.Range("K132").Value = 8082.9375
$ws.Range("L132").Value = 11136.9999
$ws.Range("M132").Value = -5552.9375
$ws.Range("N132").Value = -16196.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1536.7059
$ws.Range("I7").Value = 1536.7059
$ws.Range("K7").Value = 1536.7059
$ws.Range("M7").Value = -1424.7059
$ws.Range("H32").Value = 1281.5
$ws.Range("I32").Value = 1281.5
$ws.Range("K32").Value = 1281.5
$ws.Range("M32").Value = -964.5
$ws.Range("H40").Value = 2212.4
$ws.Range("I40").Value = 1661.08
$ws.Range("J40").Value = 4969
$ws.Range("K40").Value = 1661.08
$ws.Range("L40").Value = 4969
$ws.Range("M40").Value = -1525.08
$ws.Range("N40").Value = -5241
$ws.Range("H55").Value = 399.9524
$ws.Range("I55").Value = 243.21428
$ws.Range("J55").Value = 713.4286
$ws.Range("K55").Value = 243.21428
$ws.Range("L55").Value = 713.4286
$ws.Range("M55").Value = -70.21428
$ws.Range("N55").Value = -1059.4286
$ws.Range("H126").Value = 1536.7059
$ws.Range("I126").Value = 1536.7059
$ws.Range("K126").Value = 4610.1177
$ws.Range("M126").Value = -2140.1177
$ws.Range("H132").Value = 26902.65
$ws.Range("I132").Value = 1233.5714
$ws.Range("J132").Value = 55273.74
$ws.Range("K132").Value = 3700.7142
$ws.Range("L132").Value = 165821.22
$ws.Range("M132").Value = -1170.7142
$ws.Range("N132").Value = -170881.22
$ws.Range("H136").Value = 2378.6
$ws.Range("I136").Value = 2348.25
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 7044.75
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -4494.75
$ws.Range("N136").Value = -12600

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 701.3103599999999
$ws.Range("I136").Value = 622.86365
$ws.Range("K136").Value = 1868.59095
$ws.Range("M136").Value = 681.40905
